$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8

# Row 3 updates
$ws.Range("K3").Value = 1.73
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.18
$ws.Range("N3").Value = 4.5
$ws.Range("AC3").Value = 4.5
$ws.Range("AM3").Value = 67
$ws.Range("AT3").Value = 1.91
$ws.Range("BA3").Value = 201

# Row 5 updates
$ws.Range("G5").Value = 4.2
$ws.Range("K5").Value = 1.83
$ws.Range("Z5").Value = 41
$ws.Range("AK5").Value = 19
$ws.Range("AN5").Value = 5.5
$ws.Range("AX5").Value = 13

# Row 9 updates
$ws.Range("G9").Value = 1.98
$ws.Range("H9").Value = 3.6
$ws.Range("J9").Value = 2.57
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 3.8
$ws.Range("Q9").Value = 1.7
$ws.Range("R9").Value = 2.07
$ws.Range("W9").Value = 8.75
$ws.Range("AA9").Value = 14.5
$ws.Range("AB9").Value = 22
$ws.Range("AD9").Value = 7
$ws.Range("AE9").Value = 13
$ws.Range("AH9").Value = 12
$ws.Range("AI9").Value = 19.5
$ws.Range("AM9").Value = 30
$ws.Range("AN9").Value = 4
$ws.Range("AO9").Value = 10
$ws.Range("AU9").Value = 6.9
$ws.Range("AZ9").Value = 90
